# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on the
# zh-cn and de-de report sheets, row 2.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-02-22 04:43:38"
$zhcn.Range("G2").Value = "2016-02-22 04:44:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-02-22 04:43:52"
$dede.Range("G2").Value = "2016-02-22 04:45:09"
